{"js": "// T1499 Interworking Denial of Service \u2014 FiGHT v1.0.1 content fixes.\n// Three visible text edits:\n//  1) \"Control-plane, User-plane\" -> \"Control-plane, Roaming\" (Metadata bullet)\n//  2) \"Use WAF to minimize potential exploit of vulnerabilities\"\n//       -> \"Use Web Application Firewall (WAF) to minimize potential exploit\n//          of vulnerabilities\" with \"Web Application Firewall\" highlighted\n//          (darkCyan / teal) in the Mitigations table\n//  3) AMF description row gains a trailing clause about the N2 interface.\n\nconst body = context.document.body;\n\n// 1) Architecture Segment bullet: \"User-plane\" -> \"Roaming\"\nconst archResults = body.search(\"Control-plane, User-plane\", { matchCase: true, matchWholeWord: false });\narchResults.load(\"items\");\nawait context.sync();\nif (archResults.items.length > 0) {\n  archResults.items[0].insertText(\"Control-plane, Roaming\", Word.InsertLocation.replace);\n}\nawait context.sync();\n\n// 2) Mitigations table: spell out \"Web Application Firewall (WAF)\"\nconst wafResults = body.search(\"Use WAF to minimize potential exploit of vulnerabilities\", { matchCase: true, matchWholeWord: false });\nwafResults.load(\"items\");\nawait context.sync();\nif (wafResults.items.length > 0) {\n  wafResults.items[0].insertText(\n    \"Use Web Application Firewall (WAF) to minimize potential exploit of vulnerabilities\",\n    Word.InsertLocation.replace\n  );\n}\nawait context.sync();\n\n// Highlight the newly spelled-out phrase (teal / darkCyan)\nconst wafPhrase = body.search(\"Web Application Firewall\", { matchCase: true, matchWholeWord: false });\nwafPhrase.load(\"items\");\nawait context.sync();\nif (wafPhrase.items.length > 0) {\n  wafPhrase.items[0].font.highlightColor = \"#008080\";\n}\nawait context.sync();\n\n// 3) AMF description: append the new N2 interface clause\nconst amfResults = body.search(\"to 4G networks via N26 interface\", { matchCase: true, matchWholeWord: false });\namfResults.load(\"items\");\nawait context.sync();\nif (amfResults.items.length > 0) {\n  amfResults.items[0].insertText(\n    \"to 4G networks via N26 interface and mobility function to 5G networks via N2 interface\",\n    Word.InsertLocation.replace\n  );\n}\nawait context.sync();\n", "ps1": "# T1499 Interworking Denial of Service \u2014 FiGHT v1.0.1 content fixes.\n# Three visible text edits:\n#  1) \"Control-plane, User-plane\" -> \"Control-plane, Roaming\" (Metadata bullet)\n#  2) \"Use WAF to minimize potential exploit of vulnerabilities\"\n#       -> \"Use Web Application Firewall (WAF) to minimize potential exploit\n#          of vulnerabilities\" with \"Web Application Firewall\" highlighted\n#          (darkCyan / teal) in the Mitigations table\n#  3) AMF description row gains a trailing clause about the N2 interface.\n\n$d = $word.ActiveDocument\n\n# 1) Architecture Segment bullet: \"User-plane\" -> \"Roaming\"\n$rArch = $d.Content\n$rArch.Find.Execute(\"Control-plane, User-plane\") | Out-Null\n$rArch.Text = \"Control-plane, Roaming\"\n\n# 2) Mitigations table: spell out \"Web Application Firewall (WAF)\"\n$rWaf = $d.Content\n$rWaf.Find.Execute(\"Use WAF to minimize potential exploit of vulnerabilities\") | Out-Null\n$rWaf.Text = \"Use Web Application Firewall (WAF) to minimize potential exploit of vulnerabilities\"\n\n# Highlight the newly spelled-out phrase (teal / darkCyan, wdTeal = 10)\n$rWafHighlight = $d.Content\n$rWafHighlight.Find.Execute(\"Web Application Firewall\") | Out-Null\n$rWafHighlight.Font.HighlightColorIndex = 10\n\n# 3) AMF description: append the new N2 interface clause\n$rAmf = $d.Content\n$rAmf.Find.Execute(\"to 4G networks via N26 interface\") | Out-Null\n$rAmf.Text = \"to 4G networks via N26 interface and mobility function to 5G networks via N2 interface\"\n"}
